$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.811.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.454.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -9.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "468.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.496"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.29%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.447.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.99%  "

# Row 13
$ws.Range("E13").Value = "  -3.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.877.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.802.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.66%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.449.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -10.33%  "

# Row 19
$ws.Range("E19").Value = "  -10.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "313.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -12.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("E23").Value = "  +1.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -13.55%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "56.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.34%  "

# Row 26
$ws.Range("E26").Value = "  +1.37%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.158"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.37%  "

# Row 28
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.388"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.534.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0725"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.11%  "

# Row 35
$ws.Range("E35").Value = "  -10.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.96%  "

# Row 37
$ws.Range("E37").Value = "  -14.51%  "

# Row 38
$ws.Range("E38").Value = "  -6.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.805"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.80%  "

# Row 40
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0527"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.83%  "

# Row 44
$ws.Range("E44").Value = "  -8.57%  "

# Row 45
$ws.Range("E45").Value = "  -10.22%  "

# Row 46
$ws.Range("E46").Value = "  -2.65%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.947.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.25%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0886"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0219"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "237.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.31%  "
